# Weekly fruit/vegetable price sheet: the "Fecha" (D) + "Calidad" (L) row-pairs for
# rows 3-18 were re-sequenced (new weekly observations inserted / re-ordered while
# the underlying per-date "Primera"/"Segunda" records kept their original figures).
# Columns A, B, C, E, F, G, H, I, J, K, Q, R, T are identical for every data row, so
# only D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the columns that move, rows 3-18,
# before writing anything back out -- row N's new content may come from a
# different row N' later in the same range, so everything must be read first.
$snapshot = @{}
for ($row = 3; $row -le 18; $row++) {
    $snapshot[$row] = @{
        D = $ws.Range("D$row").Value2
        L = $ws.Range("L$row").Value2
        M = $ws.Range("M$row").Value2
        N = $ws.Range("N$row").Value2
        O = $ws.Range("O$row").Value2
        P = $ws.Range("P$row").Value2
        S = $ws.Range("S$row").Value2
    }
}

# target row -> source row (source row's original record now lives at target row)
$rowMap = @{
    3  = 17
    4  = 18
    5  = 4
    6  = 5
    7  = 15
    8  = 16
    9  = 11
    10 = 12
    11 = 10
    12 = 6
    13 = 7
    14 = 3
    15 = 13
    16 = 14
    17 = 8
    18 = 9
}

foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    $data = $snapshot[$sourceRow]

    $ws.Range("D$targetRow").Value = $data.D
    $ws.Range("L$targetRow").Value = $data.L
    $ws.Range("M$targetRow").Value = $data.M
    $ws.Range("N$targetRow").Value = $data.N
    $ws.Range("O$targetRow").Value = $data.O
    $ws.Range("P$targetRow").Value = $data.P
    $ws.Range("S$targetRow").Value = $data.S
}
